$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value for the 'Price' column (D). These must
# remain plain text cells (as in the source data), so we force a Text
# number format before assigning the value (otherwise numeric-looking
# strings like "558.35" get auto-converted to a real number by Excel),
# then clear the format again so the cell's style is left untouched.
$priceCells = [ordered]@{}
$priceCells["D2"] = "64.145.90"
$priceCells["D3"] = "3.166.24"
$priceCells["D5"] = "558.35"
$priceCells["D6"] = "171.15"
$priceCells["D9"] = "3.165.78"
$priceCells["D10"] = "0.123"
$priceCells["D11"] = "6.61"
$priceCells["D12"] = "0.395"
$priceCells["D13"] = "3.717.31"
$priceCells["D15"] = "27.49"
$priceCells["D16"] = "64.116.66"
$priceCells["D17"] = "0.0000162"
$priceCells["D18"] = "3.166.99"
$priceCells["D19"] = "5.65"
$priceCells["D20"] = "13.00"
$priceCells["D21"] = "351.57"
$priceCells["D22"] = "7.15"
$priceCells["D23"] = "1.00"
$priceCells["D24"] = "69.03"
$priceCells["D25"] = "0.500"
$priceCells["D26"] = "0.0000118"
$priceCells["D27"] = "9.40"
$priceCells["D28"] = "0.175"
$priceCells["D29"] = "1.00"
$priceCells["D30"] = "1.00"
$priceCells["D31"] = "5.58"
$priceCells["D32"] = "1.88"
$priceCells["D33"] = "22.02"
$priceCells["D34"] = "6.59"
$priceCells["D35"] = "1.19"
$priceCells["D36"] = "157.32"
$priceCells["D37"] = "1.44"
$priceCells["D38"] = "25.98"
$priceCells["D39"] = "0.795"
$priceCells["D40"] = "2.51"
$priceCells["D41"] = "1.68"
$priceCells["D42"] = "2.649.33"
$priceCells["D43"] = "6.03"
$priceCells["D44"] = "4.14"
$priceCells["D45"] = "0.0651"
$priceCells["D46"] = "38.68"
$priceCells["D47"] = "324.95"
$priceCells["D48"] = "23.51"
$priceCells["D49"] = "0.0270"
$priceCells["D51"] = "0.999"

foreach ($addr in $priceCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceCells[$addr]
    $ws.Range($addr).ClearFormats()
}

# Remaining plain text / percentage cells (safe to assign directly).
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("E3").Value = "  -7.89%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("E6").Value = "  -1.26%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -7.86%  "
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("E11").Value = "  -4.49%  "
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("E13").Value = "  -7.86%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  -4.41%  "
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").Value = "  -5.10%  "
$ws.Range("E18").Value = "  -8.04%  "
$ws.Range("E19").Value = "  -4.67%  "
$ws.Range("E20").Value = "  -6.10%  "
$ws.Range("E21").Value = "  -5.03%  "
$ws.Range("E22").Value = "  -6.58%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  -4.82%  "
$ws.Range("E25").Value = "  -6.41%  "
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("E33").Value = "  -6.81%  "
$ws.Range("E34").Value = "  -5.90%  "
$ws.Range("E35").Value = "  -7.79%  "
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("E37").Value = "  -6.02%  "
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E38").Value = "  -9.82%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E39").Value = "  -9.48%  "
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E43").Value = "  -6.75%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E44").Value = "  -6.98%  "
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("E49").Value = "  -6.97%  "
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("E51").Value = "  -0.04%  "
